$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New job-listing rows (15 rows total) for columns:
#   A = row index (0-based), B = title, C = company, D = link, E = date_listed
$titles = @(
  "newData Science Internship (Philadelphia) - Publicis Health",
  "newData Analytics Intern (Summer 2022)",
  "newSummer 2022 Data Science Intern",
  "newData Science Summer Intern - AI Innovations",
  "newResearch/Data Scientist Intern",
  "newIntern - IT Data & Analytics",
  "newHR Intern, People Analytics (Summer)",
  "newData Science Intern",
  "newInventory Planning & Business Analytics Intern",
  "newPaid Internship - IT/Tech",
  "newData Center Platform Application Engineer Intern",
  "newIOTG Research AI Scientist internship",
  "newIntern - Data Science",
  "newIntern: Energy Trading Analyst",
  "newData Science Intern"
)

$companies = @(
  "Publicis Health",
  "Poshmark",
  "Slack",
  "IBM",
  "Ascension",
  "Amcor",
  "TriNet",
  "Varian Medical Systems",
  "Fullbeauty",
  "The Shopping Center Group",
  "Intel",
  "Intel",
  "Navistar, Inc.",
  "Greenwich Commodities LLC",
  "Meketa Investment Group"
)

$links = @(
  "www.indeed.com//cmp/Publicis-Healthcare-Communications-Group",
  "www.indeed.com//cmp/Poshmark",
  "www.indeed.com//cmp/Slack",
  "www.indeed.com//cmp/IBM",
  "www.indeed.com//cmp/Ascension",
  "www.indeed.com//cmp/Amcor",
  "www.indeed.com//cmp/Trinet",
  "www.indeed.com//cmp/Varian-Medical-Systems",
  "www.indeed.com//cmp/Fullbeauty",
  "www.indeed.com//cmp/The-Shopping-Center-Group",
  "www.indeed.com//cmp/Intel-Corporation",
  "www.indeed.com//cmp/Intel-Corporation",
  "www.indeed.com//cmp/Navistar,-Inc.",
  "www.indeed.com//jobs?q=Greenwich+Commodities+LLC&l=Denver,+CO&nc=jasx",
  "www.indeed.com//cmp/Meketa-Investment-Group"
)

$dates = @(
  "PostedJust posted",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday",
  "PostedToday"
)

$startRow = 2
$n = $titles.Count

# Write column-by-column (B, then C, then D, then E) so that newly
# introduced shared strings are grouped the same way the source workbook
# groups them (all titles, then all companies, then all links, then dates).
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $titles[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $companies[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $links[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Cells.Item($startRow + $i, 5).Value = $dates[$i]
}

# Column A (row index numbers). Row 2's "A2" cell already carries the
# bold/bordered/centered template style, so copy that formatting down to
# the newly added rows before writing each row's index value.
for ($i = 0; $i -lt $n; $i++) {
    $row = $startRow + $i
    if ($row -gt $startRow) {
        $ws.Range("A$startRow").Copy($ws.Range("A$row"))
    }
    $ws.Cells.Item($row, 1).Value = $i
}
